# Fix bugs in VAR and arima: update y_value (column B) forecasts for the
# "y_fitted_on_begin_2016" and "y_fitted_on_begin_2021" sheets.

$wb = $excel.ActiveWorkbook

$sheet2016 = $wb.Worksheets.Item("y_fitted_on_begin_2016")
$sheet2021 = $wb.Worksheets.Item("y_fitted_on_begin_2021")

$values2016 = @{
    2  = -0.2964594439589135
    3  = 13.49698493759282
    4  = 12.27335075898739
    5  = 12.28369943247001
    6  = 12.57672444323161
    7  = 12.90856050812181
    8  = 12.82721742164797
    9  = 12.44595688788379
    10 = 12.00720594811706
    11 = 12.17428822689337
    12 = 11.65637214251544
    13 = 11.27365610305328
    14 = 10.69655428085121
    15 = 10.25888950036061
    16 = 10.56571049784867
    17 = 10.58569762204217
    18 = 10.03840285702023
    19 = 9.595611755396879
    20 = 8.994207046185732
    21 = 8.905739478979513
    22 = 8.157067341330361
    23 = 7.700574768518576
    24 = 7.072053664882484
    25 = 6.753995148774536
    26 = 6.311873770052049
    27 = 6.072175189988163
    28 = 6.003538030177165
}

$values2021 = @{
    2  = -0.2654151558693363
    3  = 13.52802922568239
    4  = 12.30439504707696
    5  = 12.31474372055958
    6  = 12.60776873132118
    7  = 12.93960479621138
    8  = 12.85826170973754
    9  = 12.47700117597336
    10 = 12.03825023620663
    11 = 12.20533251498294
    12 = 11.68741643060501
    13 = 11.30470039114285
    14 = 10.72759856894078
    15 = 10.28993378845018
    16 = 10.59675478593824
    17 = 10.61674191013174
    18 = 10.0694471451098
    19 = 9.626656043486456
    20 = 9.02525133427531
    21 = 8.936783767069089
    22 = 8.188111629419938
    23 = 7.731619056608155
    24 = 7.103097952972061
    25 = 6.785039436864114
    26 = 6.342918058141627
    27 = 6.10321947807774
    28 = 6.034582318266742
    29 = 5.820084173339517
    30 = 5.533322306573411
    31 = 5.511856668698503
    32 = 5.349702908406791
    33 = 5.385699053774325
}

foreach ($row in $values2016.Keys) {
    $sheet2016.Cells.Item($row, 2).Value = $values2016[$row]
}

foreach ($row in $values2021.Keys) {
    $sheet2021.Cells.Item($row, 2).Value = $values2021[$row]
}
